$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos list refresh (GitHub Actions data pull).
# Cells whose new text is numeric-looking get a forced text format so Excel
# does not silently convert them into Double values / change their display.

$ws.Range("D2").Value = '64.441.88'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '2.627.29'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.647'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.25%  '
$ws.Range("E9").Value = '  -5.09%  '
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("E11").Value = '  -2.73%  '
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.42'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.16%  '
$ws.Range("E14").Value = '  -6.14%  '
$ws.Range("D15").Value = '3.101.99'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '64.260.37'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '2.641.82'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("E19").Value = '  -2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '345.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000112'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.74'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '555.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.76%  '
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.89'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.48'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("E35").Value = '  -3.90%  '
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '154.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +4.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("E44").Value = '  -3.10%  '
$ws.Range("E45").Value = '  -2.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.38%  '
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("E48").Value = '  +3.15%  '
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0241'
$ws.Range("E50").Value = '  -6.17%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '0.0₆0241'
$ws.Range("E51").Value = '  -4.38%  '
